# Generate Report for Handback
# Update timestamps on the Overview, zh-cn, and de-de sheets to reflect
# the latest handback report generation run.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Overview!G2 - "Latest HO Xliff Generate Date" for the first data row.
$wsOverview.Range("G2").Value = "2016-09-04 05:12:01"

# de-de!H2 - "Correspond Handoff Datetime" originally held the exact same
# timestamp string as Overview!G2, so it must be refreshed in lockstep.
$wsDeDe.Range("H2").Value = "2016-09-04 05:12:01"

# zh-cn sheet: Correspond Handoff Datetime (H2) and Correspond Handback
# DateTime (K2) for the first data row.
$wsZhCn.Range("H2").Value = "2016-09-04 05:11:56"
$wsZhCn.Range("K2").Value = "2016-09-04 05:12:28"

# de-de sheet: Correspond Handback DateTime (K2) for the first data row.
$wsDeDe.Range("K2").Value = "2016-09-04 05:12:36"
